# Insert a new row above row 7, shifting existing rows 7-46 down to 8-47.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is outside the range that will be shifted, so its date format is a
# stable reference for the "Fecha" column's number format.
$dateFormat = $ws.Range("D2").NumberFormat

$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with a new price record (same shape as
# the other rows), with Fecha (D) and Volumen (J) set to the new values.
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Macroferia Regional de Talca"
$ws.Range("C7").Value = "Maule"
$ws.Range("D7").Value = 44490
$ws.Range("D7").NumberFormat = $dateFormat
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 100112022
$ws.Range("G7").Value = "Arveja Verde"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 20000
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = 20000
$ws.Range("N7").Value = "`$/saco 25 kilos"
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 800
$ws.Range("Q7").Value = 25
$ws.Range("R7").Value = "Hortaliza"
